# Applies the attendance-report sync update described in the commit:
# - Reorders the "Recorded By" email lists on several rows (content only,
#   same set of people, different order).
# - Session "ANATOMY / 5" (row 6) has now been recorded: fills in the
#   recorder, attendance count, and status, and re-colors the row from the
#   "Not Recorded" (pink) look to the "Recorded" (green) look.
# - Updates the dependent summary/statistics cells (counts, percentages)
#   to match the newly recorded session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Reordered "Recorded By" email lists (same people, new order) ----
$ws.Range("G2").Value  = "servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G3").Value  = "System, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G4").Value  = "servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G5").Value  = "eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G7").Value  = "Amera.a.saad@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg"
$ws.Range("G9").Value  = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G12").Value = "amira.m.ibrahim@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"

# ---- Row 6 (ANATOMY, session 5) moves from "Not Recorded" to "Recorded" ----
$ws.Range("G6").Value = "majorelle.magdy@med.asu.edu.eg"
$ws.Range("H6").Value = "49/251"
$ws.Range("I6").Value = "Recorded"

# Re-color A6:I6 to match the standard "Recorded" (green) row formatting,
# by copying the formatting from an existing "Recorded" row (row 2).
$ws.Range("A2:I2").Copy()
$ws.Range("A6:I6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Dependent statistics updates ----
# Subject-level stats for ANATOMY (K2:L10 block)
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 1
$ws.Range("L9").Value = "'62.1%"
$ws.Range("L10").Value = "'25.2%"

# Group statistics row (row 15)
$ws.Range("O15").Value = 18
$ws.Range("P15").Value = 1
$ws.Range("R15").Value = "'62.1%"
$ws.Range("S15").Value = "'25.2%"

# PARASITOLOGY session 4 (row 16) attendance count increased
$ws.Range("H16").Value = "114/251"
